# TC39_Canine_Filter_Breed-Rottweiler.xlsx - "updated first 25 tc in icdc breed+diagnosis"
#
# - B2 (CasesTab query): append a `Cohort` column to the RETURN clause.
# - B4 (FilesTab query): drop the trailing `Study Code` column from the RETURN clause.
# - C2/C3/C4: replace the old aliquot-count query with the new
#   Programs/Studies/Cases/Samples/Case Files/Study Files summary query.
# - Row heights shrink to fit the now-shorter cell text.
# - Column widths / selected cell nudged slightly (cosmetic, matches saved view state).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- B2: CasesTab query gains a trailing Cohort column ----
$b2 = @'
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)
WHERE demo.breed  IN ['Rottweiler']
MATCH (c)<--(diag:diagnosis)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (co:cohort)<-[*]-(c)
WITH DISTINCT c, s, demo, diag, co
RETURN  coalesce(c.case_id, '') AS `Case ID` ,
        coalesce(s.clinical_study_designation, '') AS `Study Code` ,
        coalesce(s.clinical_study_type, '') AS  `Study Type`,
        coalesce(demo.breed, '') AS Breed ,
        coalesce(diag.disease_term, '') AS Diagnosis ,
        coalesce(diag.stage_of_disease, '') AS `Stage of Disease` ,
        coalesce(demo.patient_age_at_enrollment, '') AS Age ,
        coalesce(demo.sex, '') AS Sex ,
        coalesce(demo.neutered_indicator, '') AS `Neutered Status`,
        coalesce(demo.weight, '') AS `Weight (kg)`,
        coalesce(diag.best_response, '') AS `Response to Treatment`,
coalesce(co.cohort_description, '') AS `Cohort`
'@
$ws.Range("B2").Value = $b2

# ---- B4: FilesTab query loses the trailing Study Code column ----
$b4 = @'
MATCH (f:file)-->(parent)
WITH DISTINCT f, parent
MATCH (f)-[*]->(c:case)<--(demo:demographic)
WHERE demo.breed  IN ['Rottweiler']
OPTIONAL MATCH (s:study)<-[*]-(c)<--(diag:diagnosis)
OPTIONAL MATCH (samp:sample)-->(c)
WITH DISTINCT f, parent, c, demo, diag, s
RETURN  coalesce(f.file_name, '') AS `File Name`,
        coalesce(labels(parent)[0], '') AS `Association`,
        coalesce(f.file_description, '') AS `Description`,
        coalesce(f.file_format, '') AS `Format`,
        coalesce(f.file_size, '') AS `Size`,
        coalesce(c.case_id, '') AS `Case ID`,
        coalesce(diag.disease_term,'') AS Diagnosis 
'@
$ws.Range("B4").Value = $b4

# ---- C2/C3/C4: shared "Programs/Studies/.../Study Files" summary query
#      replaces the old aliquot-count query on all three tabs ----
$cQuery = @'
MATCH (p:program)<--(s:study)<-[*]-(c:case)<--(demo:demographic)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (diag:diagnosis)-->(c)
OPTIONAL MATCH (f:file)-[*]->(c)
OPTIONAL MATCH (sf:file)-->(s)
WITH DISTINCT f, sf, samp AS samp, c, demo, diag, s, p
WHERE demo.breed  IN ['Rottweiler'] RETURN  
    count(distinct p) AS Programs,
    count(distinct s) AS Studies,
    count(distinct c) AS Cases,
    count(distinct samp) AS Samples,
    count(distinct f) AS `Case Files`,
    count(distinct sf) AS `Study Files`
'@
$ws.Range("C2").Value = $cQuery
$ws.Range("C3").Value = $cQuery
$ws.Range("C4").Value = $cQuery

# ---- Row heights: wrapped text now needs fewer lines ----
$ws.Rows.Item(2).RowHeight = 270
$ws.Rows.Item(3).RowHeight = 225
$ws.Rows.Item(4).RowHeight = 210

# ---- Column widths: tiny nudge saved along with the view state ----
$ws.Columns.Item(2).ColumnWidth = 91.66666666666667
$ws.Columns.Item(3).ColumnWidth = 74.83333333333334
$ws.Columns.Item(4).ColumnWidth = 69.5
$ws.Columns.Item(5).ColumnWidth = 27.666666666666668

# ---- Selection moves from B4 to B2 (and the frozen top-left cell resets) ----
$ws.Range("B2").Select()
